# process-course-marks-r: "some bugs fixed, function rectified."
#  1) Rename the single worksheet "Sheet 1" -> "All Marks"
#  2) Fix the C.Total column (R): it was hard-coded to 0 for every student;
#     it should be the sum of the four class-participation scores C1..C4
#     (columns N:Q).
#  3) Add a new "W.C.Total" column (weighted C.Total) right after the
#     existing "W.A.Total" column, computed as C.Total * 2.5 (rounded to
#     2 decimals, banker's/round-half-to-even rounding - matching the R
#     script that generates this workbook), and wire it into Table3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) rename sheet -------------------------------------------------
$ws.Name = "All Marks"

# --- helper: round-half-to-even (banker's rounding), like R's round() -
function RoundHalfEven([double]$value, [int]$digits) {
    $factor = [Math]::Pow(10, $digits)
    $scaled = $value * $factor
    $floorVal = [Math]::Floor($scaled)
    $diff = $scaled - $floorVal
    $epsilon = 0.0000001
    if ($diff -lt (0.5 - $epsilon)) {
        $rounded = $floorVal
    } elseif ($diff -gt (0.5 + $epsilon)) {
        $rounded = $floorVal + 1
    } else {
        if (([int64]$floorVal) % 2 -eq 0) {
            $rounded = $floorVal
        } else {
            $rounded = $floorVal + 1
        }
    }
    return $rounded / $factor
}

$firstDataRow = 2
$lastDataRow = 31

# Column numbers (1-based): N=14 O=15 P=16 Q=17 -> R=18 (C.Total)
$colC1 = 14
$colC2 = 15
$colC3 = 16
$colC4 = 17
$colCTotal = 18

# --- 2) fix C.Total (column R) = C1+C2+C3+C4 --------------------------
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $c1 = $ws.Cells.Item($r, $colC1).Value2
    $c2 = $ws.Cells.Item($r, $colC2).Value2
    $c3 = $ws.Cells.Item($r, $colC3).Value2
    $c4 = $ws.Cells.Item($r, $colC4).Value2
    $ws.Cells.Item($r, $colCTotal).Value = $c1 + $c2 + $c3 + $c4
}

# --- 3) add W.C.Total column to the table and worksheet ---------------
$lo = $ws.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()
$colWCTotal = $newCol.Index

$ws.Cells.Item(1, $colWCTotal).Value = "W.C.Total"

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $cTotal = $ws.Cells.Item($r, $colCTotal).Value2
    $weighted = RoundHalfEven ($cTotal * 2.5) 2
    $ws.Cells.Item($r, $colWCTotal).Value = $weighted
}
